$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AP (column 42) width: 23 -> 20 (stored OOXML width).
# Excel's ColumnWidth property is offset from the stored width by the
# default font's padding (5/6 of a character here), so subtract that
# offset to land exactly on a stored width of 20.
$ws.Columns.Item(42).ColumnWidth = 20 - (5/6)

# Remove all spaces from the "HD Number" values in column AP (rows 2-108),
# e.g. "HD 196761" -> "HD196761", "HD  21209A, HD  21209" -> "HD21209A,HD21209".
for ($r = 2; $r -le 108; $r++) {
    $cell = $ws.Cells.Item($r, 42)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $newVal = $val -replace " ", ""
        $cell.Value = $newVal
    }
}
